$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data values per row (rows 2-10), columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T
# (columns A-D, F, L are unchanged)
$data = @{
    2  = @{ E=3; G=27.10983833333333;  H=81.329515;  I=0.9284397459331688;  J=0.9284397459331687;  K=3; M=1.688919666666667;  N=5.066759;    O=0.02013910563427162; P=0.02013910563427162; Q=45.78633912132056;  R=412.077052091885;   S=0.01869794611840439;  T=0.01869794611840438  }
    3  = @{ E=3; G=27.10983833333333;  H=81.329515;  I=0.9284397459331688;  J=0.9284397459331687;  K=3; M=74.83790066666667;  N=224.513702;  O=0.8923860718300156;  P=0.8923860718300156;  Q=2028.843388279392;  R=18259.59049451453;  S=0.8285266978041582;   T=0.8285266978041581   }
    4  = @{ E=3; G=27.10983833333333;  H=81.329515;  I=0.9284397459331688;  J=0.9284397459331687;  K=3; M=7.335874333333334;  N=22.007623;   O=0.08747482253571279; P=0.08747482253571279; Q=198.8743672103161;  R=1789.869304892845;  S=0.08121510201060621;  T=0.08121510201060619  }
    5  = @{ E=3; G=1.190813333333333;  H=3.57244;    I=0.04078218449921273; J=0.04078218449921273; K=3; M=1.688919666666667;  N=5.066759;    O=0.02013910563427162; P=0.02013910563427162; Q=2.011188057995556;  R=18.10069252196;     S=0.0008213167216259997;T=0.0008213167216259997}
    6  = @{ E=3; G=1.190813333333333;  H=3.57244;    I=0.04078218449921273; J=0.04078218449921273; K=3; M=74.83790066666667;  N=224.513702;  O=0.8923860718300156;  P=0.8923860718300156;  Q=89.11796995254223;  R=802.06172957288;    S=0.03639345342589941;  T=0.03639345342589941  }
    7  = @{ E=3; G=1.190813333333333;  H=3.57244;    I=0.04078218449921273; J=0.04078218449921273; K=3; M=7.335874333333334;  N=22.007623;   O=0.08747482253571279; P=0.08747482253571279; Q=8.735656967791112;  R=78.62091271012001;  S=0.003567414351687331; T=0.003567414351687331 }
    8  = @{ E=3; G=0.8986996666666666; H=2.696099;   I=0.03077806956761847; J=0.03077806956761847; K=3; M=1.688919666666667;  N=5.066759;    O=0.02013910563427162; P=0.02013910563427162; Q=1.517831541460111;  R=13.660483873141;    S=0.000619842794241229; T=0.000619842794241229 }
    9  = @{ E=3; G=0.8986996666666666; H=2.696099;   I=0.03077806956761847; J=0.03077806956761847; K=3; M=74.83790066666667;  N=224.513702;  O=0.8923860718300156;  P=0.8923860718300156;  Q=67.25679638316645;  R=605.3111674484981;  S=0.027465920599958;    T=0.027465920599958    }
    10 = @{ E=3; G=0.8986996666666666; H=2.696099;   I=0.03077806956761847; J=0.03077806956761847; K=3; M=7.335874333333334;  N=22.007623;   O=0.08747482253571279; P=0.08747482253571279; Q=6.592747818075222;  R=59.334730362677;    S=0.002692306173419248; T=0.002692306173419248 }
}

foreach ($rowNum in $data.Keys) {
    $rowVals = $data[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$rowNum").Value = $rowVals[$col]
    }
}
